# Add a new worksheet "regimParam" at the end of the workbook, modeled on
# the existing "rtdata" sheet (same header/data formatting), and fill it
# with a 24-row parameter table + a CONCATENATE formula that builds a JSON
# fragment per row (same pattern already used on the other "rt*" sheets).

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("rtdata")

$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "regimParam"

# ---- formatting: copy styles from the analogous cells on rtdata ----
$src.Range("A2:C2").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)   # xlPasteFormats

$src.Range("J2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$src.Range("A3:C3").Copy()
$ws.Range("A3:C26").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- header row ----
$ws.Range("A2").Value = "paramid"
$ws.Range("B2").Value = "sensorid"
$ws.Range("C2").Value = "name"
$ws.Range("D2").Value = "["

# ---- data rows 3-26: paramid / sensorid / name ----
for ($i = 1; $i -le 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    if ($i -le 6) {
        $ws.Cells.Item($row, 2).Value = $i
    } else {
        $ws.Cells.Item($row, 2).Value = "null"
    }
    $ws.Cells.Item($row, 3).Value = "P in KC-$i"
}

# ---- formula column: D3 is the original entry, D4:D26 filled from it ----
$ws.Range("D3").Formula = '=CONCATENATE("{ ""paramid"": ",A3,",""sensorid"":",B3,", ""name"": """,C3,"""},")'
$ws.Range("D4:D26").Formula = '=CONCATENATE("{ ""paramid"": ",A4,",""sensorid"":",B4,", ""name"": """,C4,"""},")'

# ---- column widths (auto-fit to content, like Excel does on entry) ----
$ws.Columns("C:D").AutoFit()

# ---- selection state: new sheet becomes active with a prior click at E39,
# while rtdata's old selection moves to H35 and loses the tab-selected flag ----
$ws.Range("E39").Select()
$src.Activate()
$src.Range("H35").Select()
$ws.Activate()
